$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$summaryQuery = @'
SELECT DISTINCT
    COUNT(DISTINCT prg.program_id) AS "Programs",
    COUNT(DISTINCT prj.project_id) AS "Projects",
    COUNT(DISTINCT gnt.grant_id) AS "Grants",
    COUNT(DISTINCT pub.pmid) AS "Publications"
FROM 
    df_program prg
LEFT JOIN 
    df_project prj ON prg.program_id = prj."program.program_id"
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.cancer_type LIKE '%Gastric Cancer%';
'@

$programsQuery = @'
SELECT DISTINCT 
    prg.program_name AS "Program",
    prg.website AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
     CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.website       
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Gastric Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
'@

$projectsQuery = @'
SELECT DISTINCT
    prj.project_id AS "Project ID", 
    prj.project_title AS "Project Title",
    prj.org_name AS "Organization",
    prj.project_start_date AS "Project Start Date",
    prj.project_end_date AS "Project End Date"
FROM 
    df_project prj
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
     prg.cancer_type LIKE '%Gastric Cancer%'
ORDER BY 
    lower(prj.project_id) ASC
LIMIT 100;
'@

$grantsQuery = @'
SELECT DISTINCT
    gnt.grant_id AS "Grant ID", 
    prj.project_id AS "Project",
    gnt.grant_title AS "Grant Title",
    gnt.principal_investigators AS "Principal Investigators",
    gnt.program_officers AS "Program Officers",
    gnt.fiscal_year AS "Fiscal Year",
    gnt.project_end_date AS "Project End Date"
FROM 
    df_grant gnt
LEFT JOIN 
    df_project prj ON gnt."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.cancer_type LIKE '%Gastric Cancer%'
ORDER BY 
    lower(gnt.grant_id) ASC
LIMIT 100;
'@

$publicationsQuery = @'
SELECT DISTINCT
    pub.pmid AS "PubMed ID", 
    pub.title AS "Title",
    pub.authors AS "Authors",
    pub.publication_date AS "Publication Date",
    pub.cited_by AS "Cited By",
    CASE 
    WHEN pub.relative_citation_ratio = 0 THEN '0'
    WHEN pub.relative_citation_ratio = 7.0 THEN '7'
    WHEN pub.relative_citation_ratio = 2.0 THEN '2'
    WHEN pub.relative_citation_ratio = ROUND(pub.relative_citation_ratio) THEN CAST(ROUND(pub.relative_citation_ratio) AS VARCHAR) 
    ELSE CAST(ROUND(pub.relative_citation_ratio, 2) AS VARCHAR)
END AS "Relative Citation Ratio"
FROM 
    df_publication pub
LEFT JOIN 
    df_project prj ON pub."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
WHERE 
     prg.cancer_type LIKE '%Gastric Cancer%'
ORDER BY 
    lower(pub.pmid) ASC
LIMIT 100;
'@

$ws.Range("C2").Value = $summaryQuery
$ws.Range("B2").Value = $programsQuery
$ws.Range("B3").Value = $projectsQuery
$ws.Range("B4").Value = $grantsQuery
$ws.Range("B5").Value = $publicationsQuery

# Move the view: scroll so row 2 is the top-visible row, and select C3
# (mirrors the workbook's sheetView topLeftCell/selection change).
$ws.Range("C3").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
